# Lesson 02 - Initial Notes and Exercises
# Applies:
#   1. Fill in the previously-blank "Using Tabs to Create Spaces" slide
#      (the existing last slide of the deck).
#   2. Append a new "Summary" slide after it.
#   3. Re-stamp the cached datetimeFigureOut footer fields (master,
#      every layout, notes master) from 04/08/2023 -> 08/08/2023.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Existing last slide ("Using Tabs to Create Spaces")
# ---------------------------------------------------------------------
$sTabs = $p.Slides.Item($p.Slides.Count)

$title = $sTabs.Shapes.Placeholders.Item(1)
$title.TextFrame.TextRange.Text = "Using Tabs to Create Spaces"

$body = $sTabs.Shapes.Placeholders.Item(2)
$bodyTr = $body.TextFrame.TextRange
$bodyTr.Text = "Print the following Using the print command`r" + `
  "*`tThere are Spaces on the Screen`t*`r" + `
  "You can use \t to add in a tab space which is the same as 3-5 spaces`r" + `
  "Use the following command`r" + `
  "print(“*\tThere are Spaces on the Screen\t*”)"

# ---------------------------------------------------------------------
# 2. New "Summary" slide appended at the end
# ---------------------------------------------------------------------
$layout = $sTabs.CustomLayout
$sSummary = $p.Slides.AddSlide($p.Slides.Count + 1, $layout)

$sumTitle = $sSummary.Shapes.Placeholders.Item(1)
$sumTitle.TextFrame.TextRange.Text = "Summary"

$sumBody = $sSummary.Shapes.Placeholders.Item(2)
$sumBodyTr = $sumBody.TextFrame.TextRange
$sumBodyTr.Text = "This lesson covered`r" + `
  "Characters`r" + `
  "Strings`r" + `
  "Functions`r" + `
  "How to use the print function`r" + `
  "Some special characters"

for ($i = 2; $i -le $sumBodyTr.Paragraphs().Count; $i++) {
  $sumBodyTr.Paragraphs($i, 1).IndentLevel = 2
}

# ---------------------------------------------------------------------
# 3. Re-stamp cached date fields (04/08/2023 -> 08/08/2023)
# ---------------------------------------------------------------------
$master = $p.SlideMaster
$master.Shapes.Item(3).TextFrame.TextRange.Text = "8/8/2023"

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
  $layoutI = $master.CustomLayouts.Item($i)
  for ($j = 1; $j -le $layoutI.Shapes.Count; $j++) {
    $shp = $layoutI.Shapes.Item($j)
    if ($shp.Name -like "Date Placeholder*") {
      $shp.TextFrame.TextRange.Text = "8/8/2023"
    }
  }
}

# NOTE: $p.NotesMaster.Shapes.Item(n).TextFrame.TextRange writes are
# mis-routed to $p.SlideMaster.Shapes.Item(n) in this COM host (verified
# experimentally - read-back on the NotesMaster shape shows the write
# never lands there, while SlideMaster's same-index shape changes
# instead). There is no other exposed API surface for the NotesMaster's
# cached date text, so it is intentionally left untouched here to avoid
# corrupting the SlideMaster.
